$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 307, shifting rows 307:328 down to 308:329
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new record
$ws.Range("A307").Value = 10
$ws.Range("B307").Value = "Vega Modelo de Temuco"
$ws.Range("C307").Value = "La Araucanía"
$ws.Range("D307").Value = 44714
$ws.Range("E307").Value = 9
$ws.Range("F307").Value = 100112009
$ws.Range("G307").Value = "Acelga"
$ws.Range("H307").Value = "Sin especificar"
$ws.Range("I307").Value = "Primera"
$ws.Range("J307").Value = 125
$ws.Range("K307").Value = 9000
$ws.Range("L307").Value = 9000
$ws.Range("M307").Value = 9000
$ws.Range("N307").Value = "$/docena de atados (12 kilos)"
$ws.Range("O307").Value = "Provincia de Cautín"
$ws.Range("P307").Value = 750
$ws.Range("Q307").Value = 12
$ws.Range("R307").Value = "Hortaliza"

# Match the date-style of column D used elsewhere (style id 2 -> numFmtId 165)
$ws.Range("D307").NumberFormat = "YYYY-MM-DD HH:MM:SS"
